$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (columns F through M) to match new content widths.
# ColumnWidth (COM) = stored raw XML width - 5/6
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668   # F -> 26
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666   # G -> 16
$ws.Columns.Item(8).ColumnWidth = 44.166666666666664   # H -> 45
$ws.Columns.Item(9).ColumnWidth = 38.166666666666664   # I -> 39
$ws.Columns.Item(10).ColumnWidth = 37.166666666666664  # J -> 38
$ws.Columns.Item(11).ColumnWidth = 36.166666666666664  # K -> 37
$ws.Columns.Item(12).ColumnWidth = 44.166666666666664  # L -> 45
$ws.Columns.Item(13).ColumnWidth = 36.166666666666664  # M -> 37

# Row 2 - Student Leadership Council
$ws.Range("G2").Value = "(555) 345-6789"
$ws.Range("I2").Value = "https://instagram.com/studentleadersh"
$ws.Range("J2").Value = "https://facebook.com/studentleadersh"

# Row 3 - Student Ambassadors Program
$ws.Range("J3").Value = "https://facebook.com/studentambassad"

# Row 4 - Blue Ridge PASS Program
$ws.Range("F4").Value = "blueridgepasspr@brcc.edu"

# Row 5 - Social Impact Squad
$ws.Range("G5").Value = "(555) 123-4567"
$ws.Range("I5").Value = "https://instagram.com/socialimpactsqu"

# Row 6 - Collegiate FFA (CFFA)
$ws.Range("J6").Value = "https://facebook.com/collegiateffacf"
$ws.Range("L6").Value = "https://youtube.com/channel/collegiateffacf"

# Row 7 - Math Haters Club
$ws.Range("H7").Value = "https://linkedin.com/groups/mathhatersclub"

# Row 8 - Nursing Connections
$ws.Range("I8").Value = "https://instagram.com/nursingconnecti"
$ws.Range("J8").Value = "https://facebook.com/nursingconnecti"
$ws.Range("K8").Value = "https://twitter.com/nursingconnecti"
$ws.Range("M8").Value = "https://tiktok.com/@nursingconnecti"

# Row 9 - Phi Theta Kappa (PTK) International Honor Society
$ws.Range("J9").Value = "https://facebook.com/phithetakappapt"
$ws.Range("K9").Value = "https://twitter.com/phithetakappapt"

# Row 10 - STEM Club
$ws.Range("M10").Value = "https://tiktok.com/@stemclub"

# Row 11 - Veterinary Technology Club
$ws.Range("H11").Value = "https://linkedin.com/groups/veterinarytechn"
$ws.Range("K11").Value = "https://twitter.com/veterinarytechn"

# Row 12 - Adventure Club
$ws.Range("D12").Value = "https://brcc.edu/logos/adventureclub_logo.png"
$ws.Range("G12").Value = "(555) 234-5678"
$ws.Range("H12").Value = "https://linkedin.com/groups/adventureclub"
$ws.Range("J12").Value = "https://facebook.com/adventureclub"

# Row 13 - Animanga Club
$ws.Range("J13").Value = "https://facebook.com/animangaclub"

# Row 14 - Blue Ridge Christian Fellowship
$ws.Range("G14").Value = "(555) 890-1234"
$ws.Range("H14").Value = "https://linkedin.com/groups/blueridgechrist"
$ws.Range("I14").Value = "https://instagram.com/blueridgechrist"
$ws.Range("J14").Value = "https://facebook.com/blueridgechrist"
$ws.Range("M14").Value = "https://tiktok.com/@blueridgechrist"

# Row 16 - Constituting America Club
$ws.Range("D16").Value = "https://brcc.edu/logos/constitutingame_logo.png"
$ws.Range("I16").Value = "https://instagram.com/constitutingame"
$ws.Range("J16").Value = "https://facebook.com/constitutingame"

# Row 17 - Dream, Believe, Achieve (DBA) Club
$ws.Range("D17").Value = "https://brcc.edu/logos/dreambelieveach_logo.png"
$ws.Range("H17").Value = "https://linkedin.com/groups/dreambelieveach"
$ws.Range("J17").Value = "https://facebook.com/dreambelieveach"

# Row 18 - Prism Club
$ws.Range("I18").Value = "https://instagram.com/prismclub"
